$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("I3").Value = "Cmp"
$ws.Range("I3").Select() | Out-Null
